$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the editable columns as Text so numeric-looking / percent-looking
# strings are stored verbatim instead of being coerced into numbers.
# (Kept as separate contiguous ranges rather than one comma-union Range, since
# the union form here does not reliably stamp the format onto every area.)
$fmtRange1 = $ws.Range("D2:E27")
$fmtRange2 = $ws.Range("D39:E51")
$fmtRange3 = $ws.Range("B10:C11")
$fmtRange1.NumberFormat = "@"
$fmtRange2.NumberFormat = "@"
$fmtRange3.NumberFormat = "@"

$ws.Range("D2").Value = '301.73'
$ws.Range("E2").Value = '-1.17%'
$ws.Range("D3").Value = '37.46'
$ws.Range("E3").Value = '5.92%'
$ws.Range("D4").Value = '5.002'
$ws.Range("E4").Value = '-2.71%'
$ws.Range("D5").Value = '0.07850'
$ws.Range("E5").Value = '0.89%'
$ws.Range("D6").Value = '2.236'
$ws.Range("E6").Value = '-7.24%'
$ws.Range("D7").Value = '8.035'
$ws.Range("E7").Value = '0.08%'
$ws.Range("D8").Value = '4.021'
$ws.Range("E8").Value = '2.01%'
$ws.Range("D9").Value = '0.9091'
$ws.Range("E9").Value = '-1.42%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.09511'
$ws.Range("E10").Value = '-4.44%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1884'
$ws.Range("E11").Value = '4.53%'
$ws.Range("D12").Value = '0.08526'
$ws.Range("E12").Value = '-0.74%'
$ws.Range("D13").Value = '0.03523'
$ws.Range("E13").Value = '6.13%'
$ws.Range("D14").Value = '0.09954'
$ws.Range("E14").Value = '0.60%'
$ws.Range("E15").Value = '-1.04%'
$ws.Range("D16").Value = '0.005711'
$ws.Range("E16").Value = '0.81%'
$ws.Range("D17").Value = '3.467'
$ws.Range("E17").Value = '-0.11%'
$ws.Range("D18").Value = '2.074'
$ws.Range("E18").Value = '-3.15%'
$ws.Range("E19").Value = '2.91%'
$ws.Range("E20").Value = '1.14%'
$ws.Range("D21").Value = '4.774'
$ws.Range("E21").Value = '10.77%'
$ws.Range("D22").Value = '0.2203'
$ws.Range("E22").Value = '-7.55%'
$ws.Range("D23").Value = '0.04649'
$ws.Range("E23").Value = '1.74%'
$ws.Range("E24").Value = '1.03%'
$ws.Range("D25").Value = '0.004451'
$ws.Range("E25").Value = '-0.14%'
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").Value = '0.18%'
$ws.Range("D27").Value = '0.0004748'
$ws.Range("E27").Value = '28.44%'
$ws.Range("D39").Value = '0.01762'
$ws.Range("E39").Value = '-1.65%'
$ws.Range("D40").Value = '0.04751'
$ws.Range("E40").Value = '0.07%'
$ws.Range("D41").Value = '0.007833'
$ws.Range("E41").Value = '0.87%'
$ws.Range("D42").Value = '0.1392'
$ws.Range("E42").Value = '-1.30%'
$ws.Range("D43").Value = '0.007665'
$ws.Range("E43").Value = '8.01%'
$ws.Range("D44").Value = '0.002230'
$ws.Range("E44").Value = '5.35%'
$ws.Range("E45").Value = '2.97%'
$ws.Range("D46").Value = '0.00006074'
$ws.Range("E46").Value = '-0.67%'
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").Value = '0.12%'
$ws.Range("E48").Value = '217.34%'
$ws.Range("D49").Value = '0.002689'
$ws.Range("E49").Value = '34.55%'
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").Value = '0.12%'
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").Value = '0.12%'

# Drop the temporary Text number format so cells fall back to the default style
# (matches the source workbook, which carries no explicit style on these cells).
$fmtRange1.ClearFormats()
$fmtRange2.ClearFormats()
$fmtRange3.ClearFormats()
